$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Simple Taxonsorteringsordning (column B) bumps on otherwise-unchanged rows ---
$ws.Range("B2").Value  = 83225
$ws.Range("B3").Value  = 91773
$ws.Range("B4").Value  = 83225
$ws.Range("B8").Value  = 83217
$ws.Range("B9").Value  = 79002
$ws.Range("B13").Value = 79002
$ws.Range("B15").Value = 79269
$ws.Range("B20").Value = 81230
$ws.Range("B21").Value = 83091

# --- 2. Row 16 and row 17 swap their entire record content ---
# Capture the "before" values of both rows first.
$row16_A  = $ws.Range("A16").Value()
$row16_B  = $ws.Range("B16").Value()
$row16_D  = $ws.Range("D16").Value()
$row16_E  = $ws.Range("E16").Value()
$row16_F  = $ws.Range("F16").Value()
$row16_G  = $ws.Range("G16").Value()
$row16_H  = $ws.Range("H16").Value()
$row16_M  = $ws.Range("M16").Value()
$row16_P  = $ws.Range("P16").Value()
$row16_Q  = $ws.Range("Q16").Value()
$row16_R  = $ws.Range("R16").Value()
$row16_S  = $ws.Range("S16").Value()
$row16_Y  = $ws.Range("Y16").Value()
$row16_AA = $ws.Range("AA16").Value()
$row16_AC = $ws.Range("AC16").Value()
$row16_AW = $ws.Range("AW16").Value()
$row16_AX = $ws.Range("AX16").Value()

$row17_A  = $ws.Range("A17").Value()
$row17_B  = $ws.Range("B17").Value()
$row17_D  = $ws.Range("D17").Value()
$row17_E  = $ws.Range("E17").Value()
$row17_F  = $ws.Range("F17").Value()
$row17_G  = $ws.Range("G17").Value()
$row17_H  = $ws.Range("H17").Value()
$row17_P  = $ws.Range("P17").Value()
$row17_Q  = $ws.Range("Q17").Value()
$row17_R  = $ws.Range("R17").Value()
$row17_S  = $ws.Range("S17").Value()
$row17_Y  = $ws.Range("Y17").Value()
$row17_AA = $ws.Range("AA17").Value()
$row17_AW = $ws.Range("AW17").Value()
$row17_AX = $ws.Range("AX17").Value()

# Row 16 becomes what row 17 used to be (the Gammelgransskål record),
# with its taxon sort order (B) bumped by 1 (83090 -> 83091), matching the
# same-species bump already applied to row 21 above.
$ws.Range("A16").Value  = $row17_A
$ws.Range("B16").Value  = 83091
$ws.Range("D16").Value  = $row17_D
$ws.Range("E16").Value  = $row17_E
$ws.Range("F16").Value  = $row17_F
$ws.Range("G16").Value  = $row17_G
$ws.Range("H16").Value  = $row17_H
$ws.Range("K16").Value  = ""
$ws.Range("L16").Value  = ""
$ws.Range("M16").Value  = ""
$ws.Range("N16").Value  = ""
$ws.Range("P16").Value  = $row17_P
$ws.Range("Q16").Value  = $row17_Q
$ws.Range("R16").Value  = $row17_R
$ws.Range("S16").Value  = $row17_S
$ws.Range("Y16").Value  = $row17_Y
$ws.Range("AA16").Value = $row17_AA
$ws.Range("AC16").Value = ""
$ws.Range("AW16").Value = $row17_AW
$ws.Range("AX16").Value = $row17_AX

# Row 17 becomes what row 16 used to be (the Tretåig hackspett record).
# Its taxon sort order (B) is unchanged (57884 -> 57884).
$ws.Range("A17").Value  = $row16_A
$ws.Range("B17").Value  = $row16_B
$ws.Range("D17").Value  = $row16_D
$ws.Range("E17").Value  = $row16_E
$ws.Range("F17").Value  = $row16_F
$ws.Range("G17").Value  = $row16_G
$ws.Range("H17").Value  = $row16_H
$ws.Range("K17").Value  = ""
$ws.Range("L17").Value  = ""
$ws.Range("M17").Value  = $row16_M
$ws.Range("N17").Value  = ""
$ws.Range("P17").Value  = $row16_P
$ws.Range("Q17").Value  = $row16_Q
$ws.Range("R17").Value  = $row16_R
$ws.Range("S17").Value  = $row16_S
$ws.Range("Y17").Value  = $row16_Y
$ws.Range("AA17").Value = $row16_AA
$ws.Range("AC17").Value = $row16_AC
$ws.Range("AW17").Value = $row16_AW
$ws.Range("AX17").Value = $row16_AX
